$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the title placeholder shape ("ctrTitle") on the first slide
# holding the presentation's title text.
$titleShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq "Automatic Detection of Hard Exudates in Retinal Fundus Images") {
            $titleShape = $shp
            break
        }
    }
}
if ($titleShape -eq $null) {
    $titleShape = $s.Shapes.Item(1)
}

$tr = $titleShape.TextFrame.TextRange

# Re-word the title and split it across two paragraphs: the first
# paragraph keeps "Detection of Hard " and the second paragraph holds
# "Exudates in " followed by a separate run "retinal fundus images".
$tr.Text = "Detection of Hard "
[void]$tr.InsertAfter("`rExudates in ")
[void]$tr.InsertAfter("retinal fundus images")
